# edit.ps1 - Rework "Input" sheet to the new standard template column layout,
# and drop the stray empty "I" column cells on the "갑지"/"을지" summary sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Input" - reshape from the old 17-column layout (A..Q) into the
# new 16-column standard layout (A..P).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Start from a clean slate for the whole old used range (A1:Q6) so no
# leftover values from the old 17-column layout survive in cells that the
# new 16-column layout doesn't happen to overwrite (e.g. old column P).
$ws1.Range("A1:Q6").Clear()

# Protect the two date-looking text columns (A, B) from being auto-parsed
# into Excel date serials when we assign "yyyy-mm-dd" strings to .Value.
$dateRange = $ws1.Range("A2:B6")
$dateRange.NumberFormat = "@"

# New header row (A1:P1) - plain text, no bold/centered style this time.
$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws1.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Row data, columns A..O (P/"비고" is left blank for every data row).
$rows = @(
    @("2025-09-12","2025-10-14","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","5월 운반비","KS규격-1",1,0,0),
    @("2025-08-22","2025-09-27","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","5월 청구분","KS규격-2",46,4910,248446),
    @("2025-09-15","2025-08-31","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","5. 운반비","일반자재","기타","운반비","KS규격-3",1,0,0),
    @("2025-08-22","2025-10-06","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","IJ-15861","KS규격-4",1,458040,503844),
    @("2025-08-24","2025-10-16","제이비엔지니어링","제이비엔지니어링@example.com","힐스테이트 도곡동1차","delivery@example.com","힐스테이트 도곡동1차","1. 원자재비","4) ALUM. 창호","A. 압출","2월 청구분","KS규격-5",4654,5320,27235208)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    $excelRow = $r + 2
    for ($col = 0; $col -lt $rowVals.Length; $col++) {
        $ws1.Cells.Item($excelRow, $col + 1).Value = $rowVals[$col]
    }
}

# Strip leftover styling (bold header + the temporary text NumberFormat) so
# every cell ends up on the default style, matching the template. Scope each
# call tightly to the cells that actually hold data - clearing formats over
# an untouched column (e.g. P for the data rows, which stop at O) would make
# the engine materialize empty placeholder cells across that column.
$ws1.Range("A1:P1").ClearFormats()
$ws1.Range("A2:O6").ClearFormats()

# ---------------------------------------------------------------------------
# Sheets 2 & 3 ("갑지" / "을지"): remove the stray empty inline-string cells
# that were sitting in column I (rows 2-6) with no content.
# ---------------------------------------------------------------------------
for ($sheetIdx = 2; $sheetIdx -le 3; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    for ($r = 2; $r -le 6; $r++) {
        $ws.Cells.Item($r, 9).Clear()
    }
}
